# Insert a new row at row 27 (shifts existing rows 27..105 down to 28..106)
# and populate it with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").EntireRow.Insert()

$ws.Range("A27").Value = 9
$ws.Range("B27").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44608
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = 100112022
$ws.Range("G27").Value = "Arveja Verde"
$ws.Range("H27").Value = "Sin especificar"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 34
$ws.Range("K27").Value = 24000
$ws.Range("L27").Value = 26000
$ws.Range("M27").Value = 25000
$ws.Range("N27").Value = "$/saco 25 kilos"
$ws.Range("O27").Value = "Carahue"
$ws.Range("P27").Value = 1000
$ws.Range("Q27").Value = 25
$ws.Range("R27").Value = "Hortaliza"
